$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        'B' = 0.970588235294118
        'C' = 0.00150829562594268
        'D' = 0
        'E' = 0.00150829562594268
        'F' = 0.0188536953242836
        'G' = 0.000754147812971342
        'H' = 0.000754147812971342
        'I' = 0.000754147812971342
        'J' = 0.00150829562594268
        'K' = 0.00226244343891403
        'L' = 0
        'M' = 0.00150829562594268
        'N' = 0.847662141779789
        'O' = 0.655354449472097
        'P' = 0
        'Q' = 0
        'R' = 0.944947209653092
        'S' = 0.00829562594268477
        'T' = 0.0248868778280543
        'U' = 0.0165912518853695
        'V' = 0.994720965309201
        'W' = 0.960030165912519
        'X' = 0.00150829562594268
    }
    3 = @{
        'B' = 0.00226244343891403
        'C' = 0.00150829562594268
        'D' = 0.953996983408748
        'E' = 0.00754147812971342
        'F' = 0.978883861236802
        'G' = 0.0188536953242836
        'H' = 0.974358974358974
        'I' = 0.0188536953242836
        'J' = 0.000754147812971342
        'K' = 0.00678733031674208
        'L' = 0.953242835595777
        'M' = 0.00226244343891403
        'N' = 0.000754147812971342
        'O' = 0.00226244343891403
        'P' = 0.990950226244344
        'Q' = 0.00452488687782805
        'R' = 0.000754147812971342
        'S' = 0
        'T' = 0
        'U' = 0
        'V' = 0.00150829562594268
        'W' = 0.0143288084464555
        'X' = 0.990950226244344
    }
    4 = @{
        'B' = 0.00452488687782805
        'C' = 0.996983408748115
        'D' = 0.00754147812971342
        'E' = 0.0165912518853695
        'F' = 0.000754147812971342
        'G' = 0.00226244343891403
        'H' = 0.00678733031674208
        'I' = 0.00678733031674208
        'J' = 0.979638009049774
        'K' = 0.987933634992459
        'L' = 0.0248868778280543
        'M' = 0.00754147812971342
        'N' = 0.15158371040724
        'O' = 0.340874811463047
        'P' = 0
        'Q' = 0
        'R' = 0.0535444947209653
        'S' = 0.989441930618401
        'T' = 0.975113122171946
        'U' = 0.982654600301659
        'V' = 0.00150829562594268
        'W' = 0.00904977375565611
        'X' = 0.00150829562594268
    }
    5 = @{
        'B' = 0.0226244343891403
        'C' = 0
        'D' = 0.0384615384615385
        'E' = 0.974358974358974
        'F' = 0.00150829562594268
        'G' = 0.978129713423831
        'H' = 0.0180995475113122
        'I' = 0.973604826546003
        'J' = 0.0173453996983409
        'K' = 0.00301659125188537
        'L' = 0.0218702865761689
        'M' = 0.98868778280543
        'N' = 0
        'O' = 0.00150829562594268
        'P' = 0.00904977375565611
        'Q' = 0.995475113122172
        'R' = 0.000754147812971342
        'S' = 0.00150829562594268
        'T' = 0
        'U' = 0.000754147812971342
        'V' = 0.00226244343891403
        'W' = 0.0165912518853695
        'X' = 0.00603318250377074
    }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
